$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B = New-Object 'object[,]' 24,1
$arr_B[0,0] = 18.12773397212263
$arr_B[1,0] = 17.46253429964428
$arr_B[2,0] = 17.03987238317603
$arr_B[3,0] = 16.86426178491723
$arr_B[4,0] = 16.83490453684332
$arr_B[5,0] = 17.03751740433587
$arr_B[6,0] = 17.90144176139903
$arr_B[7,0] = 19.47484570427638
$arr_B[8,0] = 20.54815279214156
$arr_B[9,0] = 21.01694240224236
$arr_B[10,0] = 21.19155626360198
$arr_B[11,0] = 21.15408087550884
$arr_B[12,0] = 21.03136680703403
$arr_B[13,0] = 20.95581942689107
$arr_B[14,0] = 20.51711508049459
$arr_B[15,0] = 20.24291765413716
$arr_B[16,0] = 20.0833800132637
$arr_B[17,0] = 20.0290531332493
$arr_B[18,0] = 20.27229630598836
$arr_B[19,0] = 21.06749057342364
$arr_B[20,0] = 21.57021428766344
$arr_B[21,0] = 21.30348686851913
$arr_B[22,0] = 20.25902012062811
$arr_B[23,0] = 19.06313719162274
$ws.Range("B2:B25").Value = $arr_B

$arr_DI = New-Object 'object[,]' 24,6
$arr_DI[0,0] = 3.2826273608389
$arr_DI[0,1] = 22.74026445050044
$arr_DI[0,2] = 16.96212855333762
$arr_DI[0,3] = 17.00031330105465
$arr_DI[0,4] = 10.79167355758972
$arr_DI[0,5] = 23.47269567014709
$arr_DI[1,0] = 3.227576146732478
$arr_DI[1,1] = 22.34644934318895
$arr_DI[1,2] = 16.85748577115322
$arr_DI[1,3] = 16.70869871018251
$arr_DI[1,4] = 10.82291568123256
$arr_DI[1,5] = 23.42757438195736
$arr_DI[2,0] = 3.193550818601122
$arr_DI[2,1] = 22.10541916764979
$arr_DI[2,2] = 16.80015482332471
$arr_DI[2,3] = 16.53737018957867
$arr_DI[2,4] = 10.84483148802006
$arr_DI[2,5] = 23.40853447559392
$arr_DI[3,0] = 3.179644814283829
$arr_DI[3,1] = 22.00753605501764
$arr_DI[3,2] = 16.77855523319252
$arr_DI[3,3] = 16.46961866576537
$arr_DI[3,4] = 10.85444732271445
$arr_DI[3,5] = 23.40296659788777
$arr_DI[4,0] = 3.177333761300796
$arr_DI[4,1] = 21.99130734255343
$arr_DI[4,2] = 16.77507573342483
$arr_DI[4,3] = 16.45849687313866
$arr_DI[4,4] = 10.85608532591915
$arr_DI[4,5] = 23.40217469657773
$arr_DI[5,0] = 3.193363419807174
$arr_DI[5,1] = 22.10409751877374
$arr_DI[5,2] = 16.79985635673444
$arr_DI[5,3] = 16.5364479434814
$arr_DI[5,4] = 10.84495840024824
$arr_DI[5,5] = 23.40845049959591
$arr_DI[6,0] = 3.263701517177648
$arr_DI[6,1] = 22.60440206220609
$arr_DI[6,2] = 16.92462647230262
$arr_DI[6,3] = 16.89823763673171
$arr_DI[6,4] = 10.80187727616494
$arr_DI[6,5] = 23.4553475826614
$arr_DI[7,0] = 3.399183059282818
$arr_DI[7,1] = 23.58544736916816
$arr_DI[7,2] = 17.22300924988104
$arr_DI[7,3] = 17.66328420166752
$arr_DI[7,4] = 10.73917681869694
$arr_DI[7,5] = 23.61540694642858
$arr_DI[8,0] = 3.496353807880802
$arr_DI[8,1] = 24.29816681970491
$arr_DI[8,2] = 17.47312082396127
$arr_DI[8,3] = 18.25138592780243
$arr_DI[8,4] = 10.70651269785059
$arr_DI[8,5] = 23.77346426595658
$arr_DI[9,0] = 3.53987139015191
$arr_DI[9,1] = 24.61905159324252
$arr_DI[9,2] = 17.59315091609075
$arr_DI[9,3] = 18.52290037705358
$arr_DI[9,4] = 10.69458625285603
$arr_DI[9,5] = 23.85387966673145
$arr_DI[10,0] = 3.556238092485326
$arr_DI[10,1] = 24.73995743979077
$arr_DI[10,2] = 17.6394600195016
$arr_DI[10,3] = 18.6261516752162
$arr_DI[10,4] = 10.69049334873919
$arr_DI[10,5] = 23.88552862850948
$arr_DI[11,0] = 3.552718444761108
$arr_DI[11,1] = 24.71394711191286
$arr_DI[11,2] = 17.62944914602624
$arr_DI[11,3] = 18.60389748804796
$arr_DI[11,4] = 10.69135597418279
$arr_DI[11,5] = 23.87865965573603
$arr_DI[12,0] = 3.541220219211914
$arr_DI[12,1] = 24.62901144219163
$arr_DI[12,2] = 17.596943885295
$arr_DI[12,3] = 18.53138699007344
$arr_DI[12,4] = 10.69424103304501
$arr_DI[12,5] = 23.8564596102843
$arr_DI[13,0] = 3.534162179785222
$arr_DI[13,1] = 24.57690329348934
$arr_DI[13,2] = 17.57714368249091
$arr_DI[13,3] = 18.48702478589167
$arr_DI[13,4] = 10.69606339820772
$arr_DI[13,5] = 23.84301654098216
$arr_DI[14,0] = 3.493494831333073
$arr_DI[14,1] = 24.27711815687135
$arr_DI[14,2] = 17.4653987886754
$arr_DI[14,3] = 18.23371008905025
$arr_DI[14,4] = 10.7073512783057
$arr_DI[14,5] = 23.76837800569671
$arr_DI[15,0] = 3.468360995670151
$arr_DI[15,1] = 24.09226330744514
$arr_DI[15,2] = 17.39841869106639
$arr_DI[15,3] = 18.07923066593522
$arr_DI[15,4] = 10.71502839235325
$arr_DI[15,5] = 23.72475192034926
$arr_DI[16,0] = 3.453840903801817
$arr_DI[16,1] = 23.98563296378431
$arr_DI[16,2] = 17.36048479969826
$arr_DI[16,3] = 17.99076126673367
$arr_DI[16,4] = 10.7197200093875
$arr_DI[16,5] = 23.70046251481505
$arr_DI[17,0] = 3.448914129868809
$arr_DI[17,1] = 23.94948105025286
$arr_DI[17,2] = 17.34774386014744
$arr_DI[17,3] = 17.9608771130868
$arr_DI[17,4] = 10.72135584998138
$arr_DI[17,5] = 23.69237726955582
$arr_DI[18,0] = 3.471043238239079
$arr_DI[18,1] = 24.11197410528478
$arr_DI[18,2] = 17.40548794416821
$arr_DI[18,3] = 18.09563663009519
$arr_DI[18,4] = 10.71418257855678
$arr_DI[18,5] = 23.72931305543849
$arr_DI[19,0] = 3.5446006856766
$arr_DI[19,1] = 24.65397653739721
$arr_DI[19,2] = 17.60646858051411
$arr_DI[19,3] = 18.55267436663372
$arr_DI[19,4] = 10.69338211844551
$arr_DI[19,5] = 23.86294803023467
$arr_DI[20,0] = 3.592013307023268
$arr_DI[20,1] = 25.00461876102801
$arr_DI[20,2] = 17.74279097400772
$arr_DI[20,3] = 18.85384738744223
$arr_DI[20,4] = 10.68225625206601
$arr_DI[20,5] = 23.95725156557403
$arr_DI[21,0] = 3.566773204333757
$arr_DI[21,1] = 24.8178427987998
$arr_DI[21,2] = 17.66959309649063
$arr_DI[21,3] = 18.69292406274536
$arr_DI[21,4] = 10.68796796966205
$arr_DI[21,5] = 23.90629211349893
$arr_DI[22,0] = 3.46983081434575
$arr_DI[22,1] = 24.10306394900006
$arr_DI[22,2] = 17.4022901440662
$arr_DI[22,3] = 18.08821841507591
$arr_DI[22,4] = 10.71456410531194
$arr_DI[22,5] = 23.72724849777612
$arr_DI[23,0] = 3.362886379264974
$arr_DI[23,1] = 23.32090071170974
$arr_DI[23,2] = 17.13671737577759
$arr_DI[23,3] = 17.45117546908091
$arr_DI[23,4] = 10.75379374500879
$arr_DI[23,5] = 23.56491533269772
$ws.Range("D2:I25").Value = $arr_DI

$arr_L = New-Object 'object[,]' 24,1
$arr_L[0,0] = 10.47836564886483
$arr_L[1,0] = 10.16273145591122
$arr_L[2,0] = 9.963587186289823
$arr_L[3,0] = 9.881190863203198
$arr_L[4,0] = 9.867437009403933
$arr_L[5,0] = 9.962480856661566
$arr_L[6,0] = 10.37070164734927
$arr_L[7,0] = 11.12483632860017
$arr_L[8,0] = 11.6458591795714
$arr_L[9,0] = 11.8748619548221
$arr_L[10,0] = 11.96036600460381
$arr_L[11,0] = 11.94200610093828
$arr_L[12,0] = 11.88192109262829
$arr_L[13,0] = 11.84495742355052
$arr_L[14,0] = 11.63072660293102
$arr_L[15,0] = 11.4972042212197
$arr_L[16,0] = 11.41965469955597
$arr_L[17,0] = 11.39327078681008
$arr_L[18,0] = 11.51149610153198
$arr_L[19,0] = 11.89960293333113
$arr_L[20,0] = 12.146151068773
$arr_L[21,0] = 12.01523229013852
$arr_L[22,0] = 11.50503718503902
$arr_L[23,0] = 10.92631145802654
$ws.Range("L2:L25").Value = $arr_L

$arr_O = New-Object 'object[,]' 24,1
$arr_O[0,0] = 14.87071538585244
$arr_O[1,0] = 14.87131400086034
$arr_O[2,0] = 14.87732281440509
$arr_O[3,0] = 14.88118614909748
$arr_O[4,0] = 14.88191295793611
$arr_O[5,0] = 14.87736919475902
$arr_O[6,0] = 14.86974904281786
$arr_O[7,0] = 14.89967455299334
$arr_O[8,0] = 14.94906677997015
$arr_O[9,0] = 14.97746785345703
$arr_O[10,0] = 14.98907156414662
$arr_O[11,0] = 14.98653482823181
$arr_O[12,0] = 14.97840551521679
$arr_O[13,0] = 14.97353646939756
$arr_O[14,0] = 14.94732975018171
$arr_O[15,0] = 14.93276966244403
$arr_O[16,0] = 14.92495381049116
$arr_O[17,0] = 14.92240356263949
$arr_O[18,0] = 14.93426180837743
$arr_O[19,0] = 14.98077029531603
$arr_O[20,0] = 15.01611095428197
$arr_O[21,0] = 14.9967983240249
$arr_O[22,0] = 14.93358548029039
$arr_O[23,0] = 14.88676413336646
$ws.Range("O2:O25").Value = $arr_O
